$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 155616
$ws.Range("E2").Value = -5651
$ws.Range("F2").Value = -7378
$ws.Range("G2").Value = -10113
$ws.Range("H2").Value = -8527
$ws.Range("I2").Value = -7736
$ws.Range("J2").Value = -791
$ws.Range("K2").Value = 177358
$ws.Range("L2").Value = 157090
$ws.Range("M2").Value = 20268
$ws.Range("N2").Value = 22364
$ws.Range("O2").Value = -2096
$ws.Range("P2").Value = 9620
$ws.Range("Q2").Value = -5602
$ws.Range("R2").Value = -1992
$ws.Range("S2").Value = 5207
$ws.Range("T2").Value = 3916
$ws.Range("U2").Value = -9518
$ws.Range("V2").Value = 76490
$ws.Range("W2").Value = -3.63
$ws.Range("X2").Value = -5.48
$ws.Range("Y2").Value = -29.12
$ws.Range("Z2").Value = -4.94
$ws.Range("AA2").Value = 775.0599999999999
$ws.Range("AB2").Value = 133.34
$ws.Range("AC2").Value = -51837
$ws.Range("AD2").Value = -3.6
$ws.Range("AE2").Value = 151715
$ws.Range("AG2").Value = 1924
$ws.Range("AH2").Value = 1.03
$ws.Range("AI2").Value = -3.67
$ws.Range("AJ2").Value = 14923872
$ws.Range("D3").Value = 154436
$ws.Range("E3").Value = -21245
$ws.Range("F3").Value = -29372
$ws.Range("G3").Value = -31244
$ws.Range("H3").Value = -22092
$ws.Range("I3").Value = -20975
$ws.Range("J3").Value = -1117
$ws.Range("K3").Value = 188803
$ws.Range("L3").Value = 182615
$ws.Range("M3").Value = 6189
$ws.Range("N3").Value = 9291
$ws.Range("O3").Value = -3102
$ws.Range("P3").Value = 13721
$ws.Range("Q3").Value = -8430
$ws.Range("R3").Value = 1721
$ws.Range("S3").Value = 17729
$ws.Range("T3").Value = 1820
$ws.Range("U3").Value = -10250
$ws.Range("V3").Value = 92542
$ws.Range("W3").Value = -13.76
$ws.Range("X3").Value = -14.31
$ws.Range("Y3").Value = -132.52
$ws.Range("Z3").Value = -12.07
$ws.Range("AA3").Value = 2950.77
$ws.Range("AB3").Value = -62.48
$ws.Range("AC3").Value = -139078
$ws.Range("AD3").Value = -0.36
$ws.Range("AE3").Value = 43577
$ws.Range("AF3").Value = 1.16
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 21319818
$ws.Range("D4").Value = 128192
$ws.Range("E4").Value = -15308
$ws.Range("F4").Value = -15308
$ws.Range("G4").Value = -19617
$ws.Range("H4").Value = -27895
$ws.Range("I4").Value = -27341
$ws.Range("J4").Value = -553
$ws.Range("K4").Value = 150648
$ws.Range("L4").Value = 144055
$ws.Range("M4").Value = 6594
$ws.Range("N4").Value = 10269
$ws.Range("O4").Value = -3675
$ws.Range("P4").Value = 3329
$ws.Range("Q4").Value = -5310
$ws.Range("R4").Value = -2827
$ws.Range("S4").Value = -2009
$ws.Range("T4").Value = 1198
$ws.Range("U4").Value = -6508
$ws.Range("V4").Value = 66415
$ws.Range("W4").Value = -11.94
$ws.Range("X4").Value = -21.76
$ws.Range("Y4").Value = -279.57
$ws.Range("Z4").Value = -16.43
$ws.Range("AA4").Value = 2184.71
$ws.Range("AB4").Value = -216.63
$ws.Range("AC4").Value = -126098
$ws.Range("AD4").Value = -0.36
$ws.Range("AE4").Value = 15663
$ws.Range("AF4").Value = 2.86
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 65576960
$ws.Range("D5").Value = 111018
$ws.Range("E5").Value = 7330
$ws.Range("F5").Value = 7330
$ws.Range("G5").Value = 11890
$ws.Range("H5").Value = 6458
$ws.Range("I5").Value = 6215
$ws.Range("J5").Value = 243
$ws.Range("K5").Value = 114468
$ws.Range("L5").Value = 84561
$ws.Range("M5").Value = 29907
$ws.Range("N5").Value = 33934
$ws.Range("O5").Value = -4028
$ws.Range("P5").Value = 5383
$ws.Range("Q5").Value = -10199
$ws.Range("R5").Value = 227
$ws.Range("S5").Value = 9879
$ws.Range("T5").Value = 1088
$ws.Range("U5").Value = -11287
$ws.Range("V5").Value = 37544
$ws.Range("W5").Value = 6.6
$ws.Range("X5").Value = 5.82
$ws.Range("Y5").Value = 28.12
$ws.Range("Z5").Value = 4.87
$ws.Range("AA5").Value = 282.75
$ws.Range("AB5").Value = 51.53
$ws.Range("AC5").Value = 7447
$ws.Range("AD5").Value = 1.87
$ws.Range("AE5").Value = 31821
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 106656288
$ws.Range("D6").Value = 96444
$ws.Range("E6").Value = 10248
$ws.Range("F6").Value = 10248
$ws.Range("G6").Value = 3677
$ws.Range("H6").Value = 3201
$ws.Range("I6").Value = 3447
$ws.Range("K6").Value = 119185
$ws.Range("L6").Value = 80783
$ws.Range("M6").Value = 38402
$ws.Range("N6").Value = 38402
$ws.Range("P6").Value = 5410
$ws.Range("Q6").Value = 7251
$ws.Range("R6").Value = -132
$ws.Range("S6").Value = -5659
$ws.Range("T6").Value = 1481
$ws.Range("U6").Value = 5771
$ws.Range("V6").Value = 32053
$ws.Range("W6").Value = 10.63
$ws.Range("X6").Value = 3.32
$ws.Range("Y6").Value = 9.529999999999999
$ws.Range("Z6").Value = 2.74
$ws.Range("AA6").Value = 210.36
$ws.Range("AB6").Value = 122.94
$ws.Range("AC6").Value = 3219
$ws.Range("AD6").Value = 10.61
$ws.Range("AE6").Value = 35826
$ws.Range("AF6").Value = 0.95
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 107205752
$ws.Range("D7").Value = 82167
$ws.Range("E7").Value = 2483
$ws.Range("G7").Value = 1785
$ws.Range("H7").Value = 1646
$ws.Range("I7").Value = 1697
$ws.Range("K7").Value = 113050
$ws.Range("L7").Value = 73856
$ws.Range("M7").Value = 39193
$ws.Range("N7").Value = 39267
$ws.Range("P7").Value = 5410
$ws.Range("Q7").Value = 10972
$ws.Range("R7").Value = -476
$ws.Range("S7").Value = -3106
$ws.Range("T7").Value = 988
$ws.Range("U7").Value = 10846
$ws.Range("W7").Value = 3.02
$ws.Range("X7").Value = 2
$ws.Range("Y7").Value = 4.37
$ws.Range("Z7").Value = 1.42
$ws.Range("AA7").Value = 188.44
$ws.Range("AC7").Value = 1583
$ws.Range("AD7").Value = 15.61
$ws.Range("AE7").Value = 36633
$ws.Range("AF7").Value = 0.67
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("D8").Value = 78951
$ws.Range("E8").Value = 2025
$ws.Range("G8").Value = 1706
$ws.Range("H8").Value = 1538
$ws.Range("I8").Value = 1522
$ws.Range("K8").Value = 111403
$ws.Range("L8").Value = 70715
$ws.Range("M8").Value = 40688
$ws.Range("N8").Value = 40811
$ws.Range("P8").Value = 5410
$ws.Range("Q8").Value = 4052
$ws.Range("R8").Value = -699
$ws.Range("S8").Value = -1326
$ws.Range("T8").Value = 923
$ws.Range("U8").Value = 3075
$ws.Range("W8").Value = 2.56
$ws.Range("X8").Value = 1.95
$ws.Range("Y8").Value = 3.8
$ws.Range("Z8").Value = 1.37
$ws.Range("AA8").Value = 173.8
$ws.Range("AC8").Value = 1420
$ws.Range("AD8").Value = 17.4
$ws.Range("AE8").Value = 38074
$ws.Range("AF8").Value = 0.65
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("D9").Value = 84029
$ws.Range("E9").Value = 2675
$ws.Range("G9").Value = 2370
$ws.Range("H9").Value = 2109
$ws.Range("I9").Value = 2113
$ws.Range("K9").Value = 116832
$ws.Range("L9").Value = 74083
$ws.Range("M9").Value = 42750
$ws.Range("N9").Value = 42920
$ws.Range("P9").Value = 5410
$ws.Range("Q9").Value = 2496
$ws.Range("R9").Value = -596
$ws.Range("S9").Value = 138
$ws.Range("T9").Value = 1192
$ws.Range("U9").Value = 42
$ws.Range("W9").Value = 3.18
$ws.Range("X9").Value = 2.51
$ws.Range("Y9").Value = 5.05
$ws.Range("Z9").Value = 1.85
$ws.Range("AA9").Value = 173.29
$ws.Range("AC9").Value = 1971
$ws.Range("AD9").Value = 12.53
$ws.Range("AE9").Value = 40041
$ws.Range("AF9").Value = 0.62
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0

$ws.Range("AI6").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("AI9").ClearContents()

Write-Output "edit complete"
